# Updates "想去人数" (want-to-go headcount) figures across the four sheets,
# reflecting a refreshed data snapshot (per commit message: "output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        "F4"  = 1151
        "F10" = 628
        "F12" = 89
        "F14" = 463
        "F18" = 687
        "F19" = 2540
        "F27" = 99
        "F29" = 919
        "F31" = 64
        "F33" = 162
    }
    "演出" = @{
        "F12" = 270
        "F15" = 336
        "F16" = 336
        "F17" = 68
        "F29" = 186
    }
    "本地生活" = @{
        "F2"  = 1750
        "F5"  = 2261
        "F6"  = 899
        "F9"  = 1116
        "F11" = 73
    }
    "全部类型" = @{
        "F2"  = 1750
        "F3"  = 2261
        "F8"  = 899
        "F9"  = 1116
        "F11" = 73
        "F13" = 1151
        "F20" = 628
        "F23" = 89
        "F25" = 463
        "F28" = 687
        "F29" = 2540
        "F33" = 270
        "F35" = 99
        "F38" = 919
        "F39" = 336
        "F40" = 68
        "F42" = 64
        "F48" = 186
        "F49" = 162
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellsForSheet = $updates[$sheetName]
    foreach ($cellRef in $cellsForSheet.Keys) {
        $ws.Range($cellRef).Value = $cellsForSheet[$cellRef]
    }
}
